$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of B10 and B11 (values 3 and 4 removed), keeping their
# existing cell formatting/style intact.
$ws.Range("B10:B11").ClearContents()

# Reflect the user's final selection being on B10 when the workbook was saved.
$ws.Range("B10").Select()
